$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '65.232.15'
Set-TextValue 2 5 '  -2.62%  '

Set-TextValue 3 4 '3.669.21'
Set-TextValue 3 5 '  +3.31%  '

Set-TextValue 4 4 '0.997'
Set-TextValue 4 5 '  -0.75%  '

Set-TextValue 5 4 '407.14'
Set-TextValue 5 5 '  -2.05%  '

Set-TextValue 6 4 '134.44'
Set-TextValue 6 5 '  +4.87%  '

Set-TextValue 7 4 '3.664.50'
Set-TextValue 7 5 '  +3.55%  '

Set-TextValue 8 4 '0.622'
Set-TextValue 8 5 '  -3.52%  '

Set-TextValue 9 5 '  -0.01%  '

Set-TextValue 10 4 '0.729'
Set-TextValue 10 5 '  -5.63%  '

Set-TextValue 11 4 '0.163'
Set-TextValue 11 5 '  -7.22%  '

Set-TextValue 12 4 '0.0000319'
Set-TextValue 12 5 '  -1.78%  '

Set-TextValue 13 4 '42.34'
Set-TextValue 13 5 '  +0.29%  '

Set-TextValue 14 4 '9.98'
Set-TextValue 14 5 '  +1.50%  '

Set-TextValue 15 4 '4.235.75'
Set-TextValue 15 5 '  +2.50%  '

Set-TextValue 16 5 '  -1.40%  '

Set-TextValue 17 4 '3.685.33'
Set-TextValue 17 5 '  +1.18%  '

Set-TextValue 18 4 '20.00'
Set-TextValue 18 5 '  -0.15%  '

Set-TextValue 19 4 '13.43'
Set-TextValue 19 5 '  +9.53%  '

Set-TextValue 20 4 '1.09'
Set-TextValue 20 5 '  -2.32%  '

Set-TextValue 21 4 '65.252.83'
Set-TextValue 21 5 '  -2.52%  '

Set-TextValue 22 4 '422.64'
Set-TextValue 22 5 '  -7.95%  '

Set-TextValue 23 4 '15.29'
Set-TextValue 23 5 '  +18.91%  '

Set-TextValue 24 4 '86.16'
Set-TextValue 24 5 '  -3.47%  '

Set-TextValue 25 4 '3.01'
Set-TextValue 25 5 '  -4.16%  '

Set-TextValue 26 4 '35.97'
Set-TextValue 26 5 '  +3.82%  '

Set-TextValue 27 4 '3.20'
Set-TextValue 27 5 '  -6.95%  '

Set-TextValue 28 4 '9.43'
Set-TextValue 28 5 '  -3.58%  '

Set-TextValue 29 4 '5.14'
Set-TextValue 29 5 '  +5.61%  '

Set-TextValue 30 4 '12.68'
Set-TextValue 30 5 '  +3.98%  '

Set-TextValue 31 4 '2.72'
Set-TextValue 31 5 '  -1.37%  '

Set-TextValue 32 5 '  +2.97%  '

Set-TextValue 33 4 '6.96'
Set-TextValue 33 5 '  -3.32%  '

Set-TextValue 34 4 '0.161'
Set-TextValue 34 5 '  +3.16%  '

Set-TextValue 35 4 '41.03'
Set-TextValue 35 5 '  +5.58%  '

Set-TextValue 36 4 '55.93'
Set-TextValue 36 5 '  -0.75%  '

Set-TextValue 37 5 '  +0.21%  '

Set-TextValue 38 4 '0.0465'
Set-TextValue 38 5 '  -4.54%  '

Set-TextValue 39 4 '2.90'
Set-TextValue 39 5 '  +28.87%  '

Set-TextValue 40 2 'Stellar'
Set-TextValue 40 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 40 4 '0.140'
Set-TextValue 40 5 '  -4.22%  '

Set-TextValue 41 2 'FirstDigitalUSD'
Set-TextValue 41 3 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 41 4 '0.994'
Set-TextValue 41 5 '  -0.65%  '

Set-TextValue 42 2 'PEPE'
Set-TextValue 42 3 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 42 4 '0.0₃0650'
Set-TextValue 42 5 '  -11.42%  '

Set-TextValue 43 2 'EnergySwap'
Set-TextValue 43 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 43 4 '27.27'
Set-TextValue 43 5 '  +29.82%  '

Set-TextValue 44 2 'LidoDAOToken'
Set-TextValue 44 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 44 4 '3.33'
Set-TextValue 44 5 '  +4.71%  '

Set-TextValue 45 2 'ApeXProtocol'
Set-TextValue 45 3 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 45 4 '3.15'
Set-TextValue 45 5 '  +22.84%  '

Set-TextValue 46 2 'NEARProtocol'
Set-TextValue 46 3 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 46 4 '4.39'
Set-TextValue 46 5 '  +2.44%  '

Set-TextValue 47 4 '2.08'
Set-TextValue 47 5 '  +7.11%  '

Set-TextValue 48 2 'Monero'
Set-TextValue 48 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 48 4 '143.67'
Set-TextValue 48 5 '  -2.47%  '

Set-TextValue 49 4 '2.80'
Set-TextValue 49 5 '  -4.54%  '

Set-TextValue 50 4 '2.52'
Set-TextValue 50 5 '  -6.68%  '

Set-TextValue 51 4 '0.291'
Set-TextValue 51 5 '  -4.01%  '
